$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Add "Owyhigh Lakes from White River Road" as a new table row
$row1 = $lo.ListRows.Add()
$row1.Range.Item(1, 1).Value = "Owyhigh Lakes from White River Road"
$row1.Range.Item(1, 2).Value = 7
$row1.Range.Item(1, 3).Value = 1670
$row1.Range.Item(1, 4).Value = "moderate"

# Add "Grove of the Patriarchs and Silver Falls Loop" as a new table row
$row2 = $lo.ListRows.Add()
$row2.Range.Item(1, 1).Value = "Grove of the Patriarchs and Silver Falls Loop"
$row2.Range.Item(1, 2).Value = 5.2
$row2.Range.Item(1, 3).Value = 1000
$row2.Range.Item(1, 4).Value = "easy"

# Match the author's final selection state
$ws.Range("D30").Select()
